$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20 (row 20 was previously empty; the existing
# "Pulse time total [ms]" row, currently row 21, shifts down to row 22)
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 contents
$ws.Range("A20").Value = "Pulse time per cylinder [100th ms]"
$ws.Range("B20").Formula = "=B19*100"

# Update the selection to B9 as per the diff
$ws.Range("B9").Select()

# Adjust column A width to fit the new, longer label (best-fit recalculated by Excel)
$ws.Columns.Item(1).ColumnWidth = 28.333333333333332
